# media_tweets_categorizada.xlsx - "atualizando - tweets jornalisticos"
#
# Two new media-group rows are added to the dataset:
#   - "Quebrando o tabu" becomes the new row 2 (pushing everything else down)
#   - "Mídia NINJA" is inserted right before "Grupo Diario de Pernambuco..."
#     (originally row 11, now row 13), pushing the remaining rows down again.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert "Quebrando o tabu" as the new row 2 ---------------------------
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()

$ws.Range("A2").Value = "Quebrando o tabu"
$ws.Range("B2").Value = "esquerda; centro-esquerda"
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = 10281
$ws.Range("E2").Value = 1280
$ws.Range("F2").Value = 138
$ws.Range("G2").Value = 351
$ws.Range("H2").Value = 2056
$ws.Range("I2").Value = "Quebrando o Tabu"
$ws.Range("J2").Value = "QuebrandoOTabu"

# --- Insert "Mídia NINJA" as the new row 12 -------------------------------
# (originally row 11 "Grupo Diario de Pernambuco..." now sits at row 13
#  after the row-2 insertion above, so the new row goes in at index 12)
$ws.Rows.Item(12).Insert()
$ws.Rows.Item(12).ClearFormats()

$ws.Range("A12").Value = "Mídia NINJA"
$ws.Range("B12").Value = "esquerda; centro-esquerda"
$ws.Range("C12").Value = 3
$ws.Range("D12").Value = 638
$ws.Range("E12").Value = 177
$ws.Range("F12").Value = 12
$ws.Range("G12").Value = 12
$ws.Range("H12").Value = 213
$ws.Range("I12").Value = "Mídia NINJA"
$ws.Range("J12").Value = "MidiaNINJA"
